# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the regenerated data output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13783
$ws1.Range("F5").Value = 546
$ws1.Range("F8").Value = 1029
$ws1.Range("F9").Value = 13892
$ws1.Range("F10").Value = 14748
$ws1.Range("F23").Value = 1147
$ws1.Range("F26").Value = 5725
$ws1.Range("F27").Value = 944
$ws1.Range("F29").Value = 5411
$ws1.Range("F30").Value = 47
$ws1.Range("F31").Value = 49
$ws1.Range("F32").Value = 258

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13784
$ws4.Range("F6").Value = 546
$ws4.Range("F9").Value = 1029
$ws4.Range("F10").Value = 13892
$ws4.Range("F11").Value = 14748
$ws4.Range("F24").Value = 1147
$ws4.Range("F27").Value = 5725
$ws4.Range("F28").Value = 944
$ws4.Range("F30").Value = 5411
$ws4.Range("F31").Value = 47
$ws4.Range("F32").Value = 49
$ws4.Range("F33").Value = 258
